$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new "property_category" column before the existing "date" column
# (old H:date, I:legislator_name, J:legislator_id all shift right by one).
$ws.Range("H1").EntireColumn.Insert()

$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
$ws.Range("H3").Value = "stock"
$ws.Range("H4").Value = "stock"

# Fix stray embedded spaces in two company names
$ws.Range("B2").Value = "春源鋼鐵工業股份有限公司"
$ws.Range("B4").Value = "中華開發金融控股股份有限公司"
